{"js": "// Hotfix: update the report date on the title page from 26/05/2025 to 04/07/2025.\nconst body = context.document.body;\nconst results = body.search(\"26/05/2025\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Date '26/05/2025' not found in document body.\");\n}\n\nresults.items[0].insertText(\"04/07/2025\", \"Replace\");\nawait context.sync();\n", "ps1": "# Hotfix: update the report date on the title page from 26/05/2025 to 04/07/2025.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"26/05/2025\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"04/07/2025\",\n    1\n)\n"}
